$d = $word.ActiveDocument

# Update the date heading in the first paragraph.
$d.Paragraphs.Item(1).Range.Text = "2025-12-15 Monday"

# Update the division problems in the table, cell by cell, to avoid any
# ambiguity from duplicate text occurring in more than one cell.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "26÷5=5, 1"
$t.Cell(1, 2).Range.Text = "31÷4=7, 3"
$t.Cell(1, 3).Range.Text = "91÷5=18, 1"
$t.Cell(1, 4).Range.Text = "99÷9=11, 0"
$t.Cell(1, 5).Range.Text = "58÷7=8, 2"

$t.Cell(5, 1).Range.Text = "30÷7=4, 2"
$t.Cell(5, 2).Range.Text = "85÷2=42, 1"
$t.Cell(5, 3).Range.Text = "70÷7=10, 0"
$t.Cell(5, 4).Range.Text = "89÷7=12, 5"
$t.Cell(5, 5).Range.Text = "14÷4=3, 2"

$t.Cell(9, 1).Range.Text = "64÷3=21, 1"
$t.Cell(9, 2).Range.Text = "65÷2=32, 1"
$t.Cell(9, 3).Range.Text = "48÷9=5, 3"
$t.Cell(9, 4).Range.Text = "53÷3=17, 2"
$t.Cell(9, 5).Range.Text = "91÷4=22, 3"

$t.Cell(13, 1).Range.Text = "90÷2=45, 0"
$t.Cell(13, 2).Range.Text = "39÷3=13, 0"
$t.Cell(13, 3).Range.Text = "45÷5=9, 0"
$t.Cell(13, 4).Range.Text = "47÷6=7, 5"
$t.Cell(13, 5).Range.Text = "33÷8=4, 1"

$t.Cell(17, 1).Range.Text = "79÷7=11, 2"
$t.Cell(17, 2).Range.Text = "49÷3=16, 1"
$t.Cell(17, 3).Range.Text = "43÷6=7, 1"
$t.Cell(17, 4).Range.Text = "35÷7=5, 0"
$t.Cell(17, 5).Range.Text = "89÷2=44, 1"
